$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-12-15 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-12-16 Tuesday", 2)
$d.Content.Find.Execute("72×35=2520", $true, $false, $false, $false, $false, $true, 1, $false, "16×72=1152", 2)
$d.Content.Find.Execute("91×85=7735", $true, $false, $false, $false, $false, $true, 1, $false, "99×18=1782", 2)
$d.Content.Find.Execute("58×90=5220", $true, $false, $false, $false, $false, $true, 1, $false, "43×95=4085", 2)
$d.Content.Find.Execute("48×27=1296", $true, $false, $false, $false, $false, $true, 1, $false, "13×83=1079", 2)
$d.Content.Find.Execute("79×50=3950", $true, $false, $false, $false, $false, $true, 1, $false, "58×74=4292", 2)
$d.Content.Find.Execute("34×35=1190", $true, $false, $false, $false, $false, $true, 1, $false, "44×22=968", 2)
$d.Content.Find.Execute("60×42=2520", $true, $false, $false, $false, $false, $true, 1, $false, "65×82=5330", 2)
$d.Content.Find.Execute("42×15=630", $true, $false, $false, $false, $false, $true, 1, $false, "54×51=2754", 2)
$d.Content.Find.Execute("47×47=2209", $true, $false, $false, $false, $false, $true, 1, $false, "96×51=4896", 2)
$d.Content.Find.Execute("98×31=3038", $true, $false, $false, $false, $false, $true, 1, $false, "72×68=4896", 2)
$d.Content.Find.Execute("60×56=3360", $true, $false, $false, $false, $false, $true, 1, $false, "63×49=3087", 2)
$d.Content.Find.Execute("31×55=1705", $true, $false, $false, $false, $false, $true, 1, $false, "66×92=6072", 2)
$d.Content.Find.Execute("39×81=3159", $true, $false, $false, $false, $false, $true, 1, $false, "30×61=1830", 2)
$d.Content.Find.Execute("24×12=288", $true, $false, $false, $false, $false, $true, 1, $false, "43×76=3268", 2)
$d.Content.Find.Execute("62×52=3224", $true, $false, $false, $false, $false, $true, 1, $false, "91×66=6006", 2)
$d.Content.Find.Execute("16×50=800", $true, $false, $false, $false, $false, $true, 1, $false, "62×26=1612", 2)
$d.Content.Find.Execute("77×25=1925", $true, $false, $false, $false, $false, $true, 1, $false, "56×46=2576", 2)
$d.Content.Find.Execute("90×90=8100", $true, $false, $false, $false, $false, $true, 1, $false, "94×43=4042", 2)
$d.Content.Find.Execute("88×46=4048", $true, $false, $false, $false, $false, $true, 1, $false, "30×76=2280", 2)
$d.Content.Find.Execute("15×74=1110", $true, $false, $false, $false, $false, $true, 1, $false, "58×95=5510", 2)
$d.Content.Find.Execute("86×29=2494", $true, $false, $false, $false, $false, $true, 1, $false, "14×61=854", 2)
$d.Content.Find.Execute("32×49=1568", $true, $false, $false, $false, $false, $true, 1, $false, "87×88=7656", 2)
$d.Content.Find.Execute("53×18=954", $true, $false, $false, $false, $false, $true, 1, $false, "66×87=5742", 2)
$d.Content.Find.Execute("42×67=2814", $true, $false, $false, $false, $false, $true, 1, $false, "47×67=3149", 2)
$d.Content.Find.Execute("85×64=5440", $true, $false, $false, $false, $false, $true, 1, $false, "22×21=462", 2)
